$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Carga completa de mayo: update the month row with latest figures
$ws.Range("A2").Value = 45017
$ws.Range("B2").Value = 956
$ws.Range("C2").Value = 13
